# Populate "Sheet1" with the Jan-Jun order-items sample data, matching the
# commit's table: OrderNo / Product / Fulfilment Store / Total Price.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("A1").Value = "OrderNo"
$ws.Range("B1").Value = "Product"
$ws.Range("C1").Value = "Fulfilment Store"
$ws.Range("D1").Value = "Total Price"

# --- Row 2: Book order ---------------------------------------------------
$ws.Range("A2").Value = 41000000007
$ws.Range("B2").Value = "Book (BK123)"
$ws.Range("C2").Value = "Nowra"
$ws.Range("D2").Value = 1000
$ws.Range("D2").NumberFormat = '"$"#,##0;[Red]\-"$"#,##0'

# --- Row 3: Kettle order --------------------------------------------------
$ws.Range("A3").Value = 41000000022
$ws.Range("B3").Value = "Kettle (KT123)"
$ws.Range("C3").Value = "Blacktown"
$ws.Range("D3").Value = 20

# --- Cosmetics: column widths so the OrderNo / Product columns aren't
# truncated (mirrors the author's manual column-resize / autofit). ---------
$ws.Columns.Item(1).ColumnWidth = 11
$ws.Columns.Item(2).ColumnWidth = 14.8

# --- Leave the selection where the author's cursor ended up after typing. -
$ws.Range("D7").Select() | Out-Null
